$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Elements")

# --- Metadata sheet (sheet1) ---
# Version: 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> Alvearie Team
$ws1.Range("B9").Value = "Alvearie Team"

# Old duplicate "Contact" row (row 10) becomes "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# The second duplicate "Contact" row (row 11) is removed entirely, shifting later rows up
$ws1.Rows.Item(11).Delete()

# --- Elements sheet (sheet2) ---
# Root Extension element's Short/Definition updated to describe this specific extension
$ws2.Range("K2").Value = "Employee Termination Reason"
$ws2.Range("L2").Value = "Code indicating the reason for employee termination, if the employee was terminated (e.g., discharged, resigned)"
